$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.973.41"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.64"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.43"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("E6").Value = "  -0.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4596"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3822"
$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07728"
$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9805"
$ws.Range("E10").Value = "  +1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.24"
$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.900.16"
$ws.Range("E12").Value = "  -4.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.682"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.956"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07031"
$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.09"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009497"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.66"
$ws.Range("E19").Value = "  -1.00%  "

$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.956.85"
$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.339"
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.094"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.27"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.02"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.685"
$ws.Range("E27").Value = "  +0.66%  "

$ws.Range("E28").Value = "  +0.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.852"
$ws.Range("E29").Value = "  +1.67%  "

$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8654"
$ws.Range("E31").Value = "  +0.98%  "

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("E33").Value = "  +0.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.046"
$ws.Range("E34").Value = "  -0.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05711"
$ws.Range("E35").Value = "  +1.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.154"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("E37").Value = "  -0.70%  "

$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.037"
$ws.Range("E39").Value = "  +12.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.502"
$ws.Range("E40").Value = "  +0.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5508"
$ws.Range("E41").Value = "  -0.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000003049"
$ws.Range("E42").Value = "  +7.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1752"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.367"
$ws.Range("E44").Value = "  +1.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.233"
$ws.Range("E45").Value = "  +7.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5182"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.19"
$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06896"
$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.37"
$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  -0.63%  "
